# "raw data in appendix"
# The "No exact measurements" sheet listed a paper (VisGBT) that was raw
# conference/workshop appendix data rather than a usable entry, so the
# whole row is removed. Removing it drops the two shared strings
# ("VisGBT" and its title) that existed nowhere else in the workbook, and
# everything below shifts up by one row.

$wb = $excel.ActiveWorkbook

# Delete the VisGBT row (row 28) from the "No exact measurements" sheet.
$ws5 = $wb.Worksheets.Item("No exact measurements")
$ws5.Rows.Item(28).Delete()

# Walk every sheet and move the selection to the bottom-right cell of its
# used range (equivalent to pressing Ctrl+End), without leaving any sheet
# other than the last one marked as the active tab.
$sheetNames = @(
    "Different evaluation methodolog",
    "Different task",
    "Non-standard features",
    "No baseline performance",
    "No exact measurements",
    "Higher baseline performance",
    "No performance in chosen metric"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $null = $ws.Activate()
    $lastCell = $ws.Cells.SpecialCells(11)
    $null = $lastCell.Select()
}

$lastSheet = $wb.Worksheets.Item("No performance in chosen metric")
$null = $lastSheet.Activate()
